$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "temp solve of RWheel" - set Fitness column (C2:C12) to a constant value
$ws.Range("C2:C12").Value = 4024
